$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1834862385321101
$ws.Range("C2").Value = 0.5963302752293578
$ws.Range("J2").Value = 0.009174311926605505
$ws.Range("O2").Value = 0.004587155963302753
$ws.Range("P2").Value = 0.1146788990825688
$ws.Range("S2").Value = 0.09174311926605505
$ws.Range("B3").Value = 0.007692307692307693
$ws.Range("C3").Value = 0.02307692307692308
$ws.Range("J3").Value = 0.01538461538461539
$ws.Range("P3").Value = 0.7846153846153846
$ws.Range("S3").Value = 0.1692307692307692
$ws.Range("J4").Value = 0.025
$ws.Range("P4").Value = 0.65
$ws.Range("S4").Value = 0.325
$ws.Range("B6").Value = 0.05128205128205128
$ws.Range("D6").Value = 0.01025641025641026
$ws.Range("F6").Value = 0.1025641025641026
$ws.Range("J6").Value = 0.2307692307692308
$ws.Range("O6").Value = 0.03076923076923077
$ws.Range("Q6").Value = 0.1794871794871795
$ws.Range("R6").Value = 0.06666666666666667
$ws.Range("S6").Value = 0.3282051282051282
$ws.Range("B7").Value = 0.06666666666666667
$ws.Range("D7").Value = 0.02222222222222222
$ws.Range("F7").Value = 0.03333333333333333
$ws.Range("J7").Value = 0.1222222222222222
$ws.Range("O7").Value = 0.02222222222222222
$ws.Range("Q7").Value = 0.1888888888888889
$ws.Range("R7").Value = 0.06111111111111111
$ws.Range("S7").Value = 0.4833333333333333
$ws.Range("B8").Value = 0.07837837837837838
$ws.Range("D8").Value = 0.02702702702702703
$ws.Range("E8").Value = 0.002702702702702703
$ws.Range("F8").Value = 0.06756756756756757
$ws.Range("J8").Value = 0.08648648648648649
$ws.Range("O8").Value = 0.01351351351351351
$ws.Range("Q8").Value = 0.2108108108108108
$ws.Range("R8").Value = 0.06756756756756757
$ws.Range("S8").Value = 0.4459459459459459
$ws.Range("B9").Value = 0.08602150537634409
$ws.Range("F9").Value = 0.07526881720430108
$ws.Range("J9").Value = 0.1075268817204301
$ws.Range("O9").Value = 0.03763440860215054
$ws.Range("Q9").Value = 0.1774193548387097
$ws.Range("R9").Value = 0.1129032258064516
$ws.Range("S9").Value = 0.4032258064516129
$ws.Range("B10").Value = 0.1190738699007718
$ws.Range("D10").Value = 0.02756339581036384
$ws.Range("F10").Value = 0.07497243660418963
$ws.Range("J10").Value = 0.09040793825799338
$ws.Range("O10").Value = 0.009922822491730982
$ws.Range("Q10").Value = 0.1984564498346196
$ws.Range("R10").Value = 0.08269018743109151
$ws.Range("S10").Value = 0.3969128996692393
$ws.Range("G11").Value = 0.1642335766423358
$ws.Range("J11").Value = 0.06569343065693431
$ws.Range("K11").Value = 0.2007299270072993
$ws.Range("L11").Value = 0.5474452554744526
$ws.Range("S11").Value = 0.0218978102189781
$ws.Range("G12").Value = 0.7161290322580646
$ws.Range("J12").Value = 0.2
$ws.Range("K12").Value = 0.01935483870967742
$ws.Range("L12").Value = 0.02580645161290323
$ws.Range("S12").Value = 0.03870967741935484
$ws.Range("G13").Value = 0.7428571428571429
$ws.Range("J13").Value = 0.2571428571428571
$ws.Range("F15").Value = 0.01595744680851064
$ws.Range("H15").Value = 0.148936170212766
$ws.Range("I15").Value = 0.0851063829787234
$ws.Range("J15").Value = 0.351063829787234
$ws.Range("K15").Value = 0.09042553191489362
$ws.Range("M15").Value = 0.01595744680851064
$ws.Range("O15").Value = 0.09574468085106383
$ws.Range("S15").Value = 0.1968085106382979
$ws.Range("F16").Value = 0.02027027027027027
$ws.Range("H16").Value = 0.1891891891891892
$ws.Range("I16").Value = 0.06081081081081081
$ws.Range("J16").Value = 0.3851351351351351
$ws.Range("K16").Value = 0.1216216216216216
$ws.Range("M16").Value = 0.03378378378378379
$ws.Range("O16").Value = 0.08108108108108109
$ws.Range("S16").Value = 0.1081081081081081
$ws.Range("F17").Value = 0.02506963788300836
$ws.Range("H17").Value = 0.181058495821727
$ws.Range("I17").Value = 0.1030640668523677
$ws.Range("J17").Value = 0.3983286908077994
$ws.Range("K17").Value = 0.09192200557103064
$ws.Range("M17").Value = 0.01114206128133705
$ws.Range("O17").Value = 0.06685236768802229
$ws.Range("S17").Value = 0.1225626740947075
$ws.Range("F18").Value = 0.01388888888888889
$ws.Range("H18").Value = 0.2013888888888889
$ws.Range("I18").Value = 0.1388888888888889
$ws.Range("J18").Value = 0.3263888888888889
$ws.Range("K18").Value = 0.0763888888888889
$ws.Range("M18").Value = 0.01388888888888889
$ws.Range("N18").Value = 0.01388888888888889
$ws.Range("O18").Value = 0.0763888888888889
$ws.Range("S18").Value = 0.1388888888888889
$ws.Range("F19").Value = 0.01890359168241966
$ws.Range("H19").Value = 0.2098298676748582
$ws.Range("I19").Value = 0.09924385633270322
$ws.Range("J19").Value = 0.3223062381852552
$ws.Range("K19").Value = 0.1247637051039698
$ws.Range("M19").Value = 0.01984877126654064
$ws.Range("N19").Value = 0.000945179584120983
$ws.Range("O19").Value = 0.06899810964083176
$ws.Range("S19").Value = 0.1351606805293006

Write-Host "Applied all cell updates"
